# Circle Language Spec Plan: Set font to Calibri for non-heading text.
# Also relocate the "_GoBack" bookmark from the title to the end of the
# "2008-03-06 - 2008-03-09" date-range paragraph (this is what Word does
# automatically: it marks the location of the most recent edit).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark.
#    It currently sits between "2008-03 " and "Computer Language..." in
#    the title. Delete it there, then re-create it right after the
#    "2008-03-06 - 2008-03-09" text (last edited location).
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$findRng = $d.Content
$found = $findRng.Find.Execute("2008-03-06 – 2008-03-09", $true, $true)

if ($found) {
    $targetStart = $findRng.Start
    $targetEnd = $findRng.End

    # Build a fresh Range object (re-using the Find range directly confuses
    # InsertXML and duplicates text), then re-insert its own text together
    # with the bookmark markers appended right after it, so the bookmark
    # ends up as a zero-width point right after the text, before </w:p>.
    $targetRng = $d.Range($targetStart, $targetEnd)
    $targetText = $targetRng.Text
    # Escape any XML-significant characters in case of future edits.
    $targetText = $targetText.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")

    $bookmarkXml = '<?xml version="1.0"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:r><w:t>' + $targetText + '</w:t></w:r>' +
        '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $targetRng.InsertXML($bookmarkXml)
}

# ---------------------------------------------------------------------
# 2) Normal style: switch the document's base font from Tahoma 12pt to
#    Calibri 11pt (non-heading text).
# ---------------------------------------------------------------------
$normalStyle = $d.Styles("Normal")
$normalStyle.Font.Name = "Calibri"
$normalStyle.Font.Size = 11
